$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three "author / date" stamp cells (multi-line: name + date,
# separated by an in-cell line break just like the original entries).
$ws.Range("G2").Value = "Dev Soni" + [char]10 + "Aug 8th, 2022"
$ws.Range("G12").Value = "Dev Soni" + [char]10 + "Aug 12, 2022"
$ws.Range("G22").Value = "Dev Soni" + [char]10 + "Aug 17, 2022"

# Leave the view scrolled/selected on the last-edited cell, as in the
# original edit (selection ends up as a single cell, G22).
$ws.Range("G22").Select()
